# Fruta / hortaliza, semanal
# Insert two new weekly price rows for Kiwi (Hayward) at "Provincia de Curicó"
# ahead of the existing block, shifting the previous 8 rows (468-475) down to
# (470-477).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 468.
$ws.Rows("468:469").Insert()

# New row 468
$ws.Range("A468").Value = 9
$ws.Range("B468").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C468").Value = "Metropolitana"
$ws.Range("D468").Value = 44656
$ws.Range("E468").Value = 13
$ws.Range("F468").Value = "Fruta"
$ws.Range("G468").Value = 100101
$ws.Range("H468").Value = "Berries"
$ws.Range("I468").Value = 100101007
$ws.Range("J468").Value = "Kiwi"
$ws.Range("K468").Value = "Hayward"
$ws.Range("L468").Value = "Primera"
$ws.Range("M468").Value = 300
$ws.Range("N468").Value = 18000
$ws.Range("O468").Value = 18000
$ws.Range("P468").Value = 18000
$ws.Range("Q468").Value = "`$/caja 18 kilos"
$ws.Range("R468").Value = "Provincia de Curicó"
$ws.Range("S468").Value = 1000
$ws.Range("T468").Value = 18

# New row 469
$ws.Range("A469").Value = 9
$ws.Range("B469").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C469").Value = "Metropolitana"
$ws.Range("D469").Value = 44656
$ws.Range("E469").Value = 13
$ws.Range("F469").Value = "Fruta"
$ws.Range("G469").Value = 100101
$ws.Range("H469").Value = "Berries"
$ws.Range("I469").Value = 100101007
$ws.Range("J469").Value = "Kiwi"
$ws.Range("K469").Value = "Hayward"
$ws.Range("L469").Value = "Segunda"
$ws.Range("M469").Value = 280
$ws.Range("N469").Value = 12000
$ws.Range("O469").Value = 12000
$ws.Range("P469").Value = 12000
$ws.Range("Q469").Value = "`$/caja 18 kilos"
$ws.Range("R469").Value = "Provincia de Curicó"
$ws.Range("S469").Value = 667
$ws.Range("T469").Value = 18

# Ensure date cells keep the same date/time number format as the rest of
# column D (style index carried over from the source row by Insert, but set
# explicitly to be safe).
$ws.Range("D468:D469").NumberFormat = $ws.Range("D470").NumberFormat
